$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (invitation 20349, guest "Guest"/"אורח/ת"): RSVP action reset back to default ("Yes")
# and diet info was added.
$ws.Range("E6").Value = "Yes"
$ws.Range("I6").Value = "Vegetarian"
